$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/date-like updates (unambiguous strings - no special handling needed)
$ws.Range("H2").Value = "2019-09-30 00:00:00"
$ws.Range("AC2").Value = "2019Q3"
$ws.Range("AD2").Value = "2019年 三季报"

# Pure numeric cell updates
$ws.Range("I2").Value = 0.12
$ws.Range("K2").Value = 44576021.69
$ws.Range("L2").Value = 3588039.88
$ws.Range("R2").Value = 36.4240691619

# Cells that must hold a NUMERIC-LOOKING string as TEXT (not a number).
# A leading apostrophe forces Excel to store it as text, then we re-apply
# the style of an already-plain-text cell (J2) so the quote-prefix
# formatting flag doesn't leave a stray style behind.
$ws.Range("AB2").Value = "'0"
$ws.Range("AB2").Style = $ws.Range("J2").Style

$ws.Range("AE2").Value = "'2019"
$ws.Range("AE2").Style = $ws.Range("J2").Style

# N2, O2, P2, Q2 become empty text cells (previously numeric).
$ws.Range("N2").Value = "'"
$ws.Range("N2").Style = $ws.Range("J2").Style

$ws.Range("O2").Value = "'"
$ws.Range("O2").Style = $ws.Range("J2").Style

$ws.Range("P2").Value = "'"
$ws.Range("P2").Style = $ws.Range("J2").Style

$ws.Range("Q2").Value = "'"
$ws.Range("Q2").Style = $ws.Range("J2").Style
